# Update the "Forecast Comparison" sheet with the correct forecast output.
# - Insert a new "Week_Start_Date" column right after "Week" (shifts ASIN..is_holiday_week
#   one column to the right: B->C, C->D, D->E, E->F, F->G, G->H, H->I, I->J).
# - Normalize the Week labels from "W01".."W16" to "W1".."W16".
# - Populate the new Week_Start_Date column with the corresponding ISO week-start dates.
# - Convert the is_holiday_week column values to real boolean FALSE values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert a new column B (Week_Start_Date), shifting everything after "Week" right ---
$ws.Columns.Item(2).Insert()

# --- Header row ---
$ws.Range("B1").Value = "Week_Start_Date"

# --- Week label + Week start date per row (rows 2..17) ---
$weekData = @(
    @("W1", "2025-01-05"),
    @("W2", "2025-01-12"),
    @("W3", "2025-01-19"),
    @("W4", "2025-01-26"),
    @("W5", "2025-02-02"),
    @("W6", "2025-02-09"),
    @("W7", "2025-02-16"),
    @("W8", "2025-02-23"),
    @("W9", "2025-03-02"),
    @("W10", "2025-03-09"),
    @("W11", "2025-03-16"),
    @("W12", "2025-03-23"),
    @("W13", "2025-03-30"),
    @("W14", "2025-04-06"),
    @("W15", "2025-04-13"),
    @("W16", "2025-04-20")
)

# Force column B to Text format first so the ISO date strings are stored as plain
# text (e.g. "2025-01-05") instead of being auto-converted to date serial numbers.
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt $weekData.Count; $i++) {
    $rowNum = $i + 2
    $week = $weekData[$i][0]
    $weekStart = $weekData[$i][1]

    $ws.Cells.Item($rowNum, 1).Value = $week        # A: Week  ("W01" -> "W1", etc.)
    $ws.Cells.Item($rowNum, 2).Value = $weekStart    # B: Week_Start_Date

    # J: is_holiday_week -> real boolean FALSE
    $ws.Cells.Item($rowNum, 10).Value = $false
}

$wb.Save()
